# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 91 (pushing the previous rows 91-120
# down to 92-121) and populate it with the new Sandia / "Extra" quality
# record for "Región de O'Higgins".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 91, shifting existing data
# (rows 91-120) down to rows 92-121.
$ws.Rows.Item(91).Insert()

# Populate the newly inserted row 91 with the new record's values.
$ws.Cells.Item(91, 1).Value = 8
$ws.Cells.Item(91, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(91, 3).Value = "Coquimbo"
$ws.Cells.Item(91, 4).Value = 44917
$ws.Cells.Item(91, 5).Value = 4
$ws.Cells.Item(91, 6).Value = 100112028
$ws.Cells.Item(91, 7).Value = "Sandia"
$ws.Cells.Item(91, 8).Value = "Sin especificar"
$ws.Cells.Item(91, 9).Value = "Extra"
$ws.Cells.Item(91, 10).Value = 1200
$ws.Cells.Item(91, 11).Value = 4000
$ws.Cells.Item(91, 12).Value = 4500
$ws.Cells.Item(91, 13).Value = 4250
$ws.Cells.Item(91, 14).Value = "$/unidad"
$ws.Cells.Item(91, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(91, 16).Value = 4250
$ws.Cells.Item(91, 17).Value = 1
$ws.Cells.Item(91, 18).Value = "Hortaliza"
